$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the leaving_time for the existing last row (row 436): 15:30 -> 16:30
$ws.Range("D436").Value = "16:30"

# 2) Append new attendance rows 437-444
#    Columns: A=staff_name, B=id, C=start_time, D=leaving_time, E=reason,
#             F=department, G=date, H=btn_id
$newRows = @(
    @{ Row=437; A="LÊ MINH THẮNG";        B="223906"; C="16:30"; D="16:30"; E="TEST REQUEST"; F="RD"; G="2024-05-18"; H="btn_12" },
    @{ Row=438; A="LÊ QUỐC TRUNG";        B="224016"; C="16:30"; D="19:0";  E="B/T";          F="RD"; G="2024-05-18"; H="btn_13" },
    @{ Row=439; A="NGUYỄN HOÀNG VIỆT";    B="172759"; C="16:30"; D="19:30"; E="TEST REQUEST"; F="RD"; G="2024-05-18"; H="btn_6"  },
    @{ Row=440; A="TRẦN VĂN LƯU";         B="234102"; C="16:30"; D="19:30"; E="TEST REQUEST"; F="RD"; G="2024-05-18"; H="btn_15" },
    @{ Row=441; A="PHẠM THỊ PHƯƠNG";      B="172684"; C="16:30"; D="19:30"; E="TEST REQUEST"; F="RD"; G="2024-05-18"; H="btn_5"  },
    @{ Row=442; A="NGUYỄN QUANG QUÍ";     B="203638"; C="16:30"; D="19:30"; E="TEST REQUEST"; F="RD"; G="2024-05-18"; H="btn_10" },
    @{ Row=443; A="LÊ MINH THẮNG";        B="223906"; C="7:30";  D="16:30"; E="TEST REQUEST"; F="RD"; G="2024-05-19"; H="btn_12" },
    @{ Row=444; A="LÊ HUỲNH ANH KHOA";    B="234168"; C="7:30";  D="16:30"; E="B/T";          F="RD"; G="2024-05-19"; H="btn_16" }
)

foreach ($r in $newRows) {
    $rowRange = $ws.Range("A$($r.Row):H$($r.Row)")
    # Force text format first so id-like and date-like strings are not
    # auto-converted to numbers / dates by Excel's smart entry parsing.
    $rowRange.NumberFormat = "@"

    $ws.Range("A$($r.Row)").Value = $r.A
    $ws.Range("B$($r.Row)").Value = $r.B
    $ws.Range("C$($r.Row)").Value = $r.C
    $ws.Range("D$($r.Row)").Value = $r.D
    $ws.Range("E$($r.Row)").Value = $r.E
    $ws.Range("F$($r.Row)").Value = $r.F
    $ws.Range("G$($r.Row)").Value = $r.G
    $ws.Range("H$($r.Row)").Value = $r.H

    # Restore the default (unstyled) cell appearance to match the rest of
    # the sheet's plain data rows, now that the text values are locked in.
    $rowRange.Style = "Normal"
}
